$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "Through 2021-09-09"
$ws.Range("B1").Value = "September 2021 (through September 09)"

$ws.Range("B2").Value = 6
$ws.Range("AL6").Value = 1
$ws.Range("AL12").Value = 1
$ws.Range("AC17").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("AU18").Value = 1
$ws.Range("K23").Value = 1
$ws.Range("AU24").Value = 2
$ws.Range("B28").Value = 2
$ws.Range("AL29").Value = 1
$ws.Range("T30").Value = 1
$ws.Range("B42").Value = 1
$ws.Range("B44").Value = 1
$ws.Range("AC44").Value = 1
$ws.Range("B73").Value = 1
$ws.Range("AC91").Value = 1
